$wb = $excel.ActiveWorkbook

$ws0 = $wb.Worksheets.Item("Layer0")
$ws0.Range("B2").Value = -0.6191948357417691
$ws0.Range("C2").Value = 0.4159850722554189
$ws0.Range("B3").Value = 0.6106777816265282
$ws0.Range("C3").Value = -1.244048072573338
$ws0.Range("B4").Value = -1.036296691387035
$ws0.Range("C4").Value = -1.335196263520278

$ws1 = $wb.Worksheets.Item("Layer1")
$ws1.Range("B2").Value = -0.4058235718665035
$ws1.Range("C2").Value = -0.01019866353021514
$ws1.Range("B3").Value = -0.2581535674892162
$ws1.Range("C3").Value = 0.5716470043357291
$ws1.Range("B4").Value = -2.134410621890477
$ws1.Range("C4").Value = -0.5780563044708875
